# Deep Audit complete: Removed all hardcoded defaults
# Zero out pre-filled/default sample values across the Purchasing Dashboard.

$wb = $excel.ActiveWorkbook

# --- SUPPLIER_CONFIG sheet: TABLE 1 (Parts Suppliers, rows 6-11, cols C:F) ---
$wsSupplier = $wb.Worksheets.Item("SUPPLIER_CONFIG")

$wsSupplier.Range("C6:F11").Value = 0

# --- SUPPLIER_CONFIG sheet: TABLE 2 (Pieces Configuration, rows 16-21, cols B:C) ---
$wsSupplier.Range("B16:C21").Value = 0

# --- MRP_ENGINE sheet: Target Production row (row 6, cols B:I) ---
$wsMrp = $wb.Worksheets.Item("MRP_ENGINE")

$wsMrp.Range("B6:I6").Value = 0

# --- MRP_ENGINE sheet: Order Supplier labels (Lead/Batch defaults reset to 0) ---
$wsMrp.Range("A27").Value = "Order Supplier A (Lead:0, Batch:0)"
$wsMrp.Range("A28").Value = "Order Supplier B (Lead:0, Batch:0)"
$wsMrp.Range("A29").Value = "Order Supplier C (Lead:0, Batch:0)"
$wsMrp.Range("A33").Value = "Order Supplier A (Lead:0, Batch:0)"
$wsMrp.Range("A34").Value = "Order Supplier B (Lead:0, Batch:0)"
$wsMrp.Range("A35").Value = "Order Supplier C (Lead:0, Batch:0)"
